$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1) Insert a new "2022-Q1" worksheet right before the "总计" (total) sheet
# ----------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2021-Q2")

$newSheet = $wb.Worksheets.Add($total)
$newSheet.Name = "2022-Q1"

# Copy header row (B1:H1) and data-row index column (A2:A3) formatting from
# an existing quarter sheet so the new sheet matches the established style.
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$template.Range("A2:A3").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns B,D,E,F,G hold numeric-looking text (fund codes / percentages
# stored as strings with significant leading/trailing digits), so force
# text format before assignment to keep them as text instead of Excel
# auto-converting them to numbers.
$newSheet.Range("B2:B3").NumberFormat = "@"
$newSheet.Range("D2:G3").NumberFormat = "@"

# Row 2
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "005607"
$newSheet.Range("C2").Value = "华宝中证500指数增强A"
$newSheet.Range("D2").Value = "0.45"
$newSheet.Range("E2").Value = "94.72"
$newSheet.Range("F2").Value = "1.60"
$newSheet.Range("G2").Value = "0.0072"
$newSheet.Range("H2").Value = 7

# Row 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "005608"
$newSheet.Range("C3").Value = "华宝中证500指数增强C"
$newSheet.Range("D3").Value = "0.23"
$newSheet.Range("E3").Value = "94.72"
$newSheet.Range("F3").Value = "1.60"
$newSheet.Range("G3").Value = "0.0037"
$newSheet.Range("H3").Value = 7

# ----------------------------------------------------------------------
# 2) Insert a new top data row in "总计" for 2022-Q1, shifting the
#    existing rows down, and renumber the index column A.
# ----------------------------------------------------------------------
# Re-fetch the "总计" sheet by name: after Worksheets.Add() the earlier
# $total handle now resolves to the newly inserted "2022-Q1" sheet
# (it was positional), so grab a fresh reference here.
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert(-4121)
$total.Range("B2:D2").ClearFormats()

# A2 picked up no style from the insert; copy it from the row below (A3),
# which still carries the original index-column style.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.01

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
